# Add data for 2021-11-06
# Updates the "through" date in the sheet name and October row label,
# and refreshes the August/September/October/Total figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet to reflect new "through" date
$wb.Sheets.Item(1).Name = "Through 2021-10-29"

# Row 10 - August 2021 (column U/V = 2021 arrest_made/no_arrest_made/arrest_rate group)
$ws.Range("U10").Value = 152
$ws.Range("V10").Value = 0.05

# Row 11 - September
$ws.Range("T11").Value = 6
$ws.Range("U11").Value = 172
$ws.Range("V11").Value = 0.0337

# Row 12 - October (through 10-29), update label and all year columns
$ws.Range("A12").Value = "October (through 10-29)"

$ws.Range("C12").Value = 27
$ws.Range("D12").Value = 0.06900000000000001

$ws.Range("E12").Value = 6
$ws.Range("F12").Value = 47
$ws.Range("G12").Value = 0.1132

$ws.Range("I12").Value = 67
$ws.Range("J12").Value = 0.141

$ws.Range("L12").Value = 55
$ws.Range("M12").Value = 0.0833

$ws.Range("O12").Value = 53
$ws.Range("P12").Value = 0.0702

$ws.Range("R12").Value = 138
$ws.Range("S12").Value = 0.0072

$ws.Range("T12").Value = 1
$ws.Range("U12").Value = 180
$ws.Range("V12").Value = 0.0055
$ws.Range("V12").NumberFormat = $ws.Range("S12").NumberFormat

# Row 13 - Total
$ws.Range("C13").Value = 223
$ws.Range("D13").Value = 0.1255

$ws.Range("E13").Value = 52
$ws.Range("F13").Value = 430
$ws.Range("G13").Value = 0.1079

$ws.Range("I13").Value = 644
$ws.Range("J13").Value = 0.08649999999999999

$ws.Range("L13").Value = 542
$ws.Range("M13").Value = 0.1086

$ws.Range("O13").Value = 432
$ws.Range("P13").Value = 0.09810000000000001

$ws.Range("R13").Value = 986
$ws.Range("S13").Value = 0.0519

$ws.Range("T13").Value = 84
$ws.Range("U13").Value = 1345
$ws.Range("V13").Value = 0.0588
